$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: columns B, C, E hold non-numeric text - set directly.
# Column D holds numeric-looking values that must stay as TEXT (matching the
# source data, which stores prices as inline strings, not numbers). Force text
# by setting NumberFormat to "@" before assigning, then ClearFormats() to avoid
# leaving a stray number-format style behind (matches original unstyled cells).

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '244.16'
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '22.37'
$cell.ClearFormats()
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '5.433'
$cell.ClearFormats()
$ws.Range("E4").Value = '3HuobiTokenHTBestin24h'
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.05762'
$cell.ClearFormats()
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '3.429'
$cell.ClearFormats()
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '6.339'
$cell.ClearFormats()
$ws.Range("E7").Value = '6KuCoinTokenKCS'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.8136'
$cell.ClearFormats()
$ws.Range("E8").Value = '7MXTokenMX'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.8758'
$cell.ClearFormats()
$ws.Range("E9").Value = '8FTXTokenFTT'
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07354'
$cell.ClearFormats()
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.03055'
$cell.ClearFormats()
$ws.Range("B13").Value = 'ProBitToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.1321'
$cell.ClearFormats()
$ws.Range("E13").Value = '12ProBitTokenPROB'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.03101'
$cell.ClearFormats()
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.09409'
$cell.ClearFormats()
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '3.932'
$cell.ClearFormats()
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.001575'
$cell.ClearFormats()
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.04828'
$cell.ClearFormats()
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.0005844'
$cell.ClearFormats()
$ws.Range("E19").Value = '18OneONE'
$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.006328'
$cell.ClearFormats()
$ws.Range("E20").Value = '19TigerCashTCH'
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.004136'
$cell.ClearFormats()
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("B22").Value = 'BitKan'
$ws.Range("C22").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.0009985'
$cell.ClearFormats()
$ws.Range("E22").Value = '21BitKanKAN'
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.0001501'
$cell.ClearFormats()
$ws.Range("E23").Value = '22NitroExNTX'
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '3.724'
$cell.ClearFormats()
$ws.Range("E24").Value = '23LEOLEO'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.196'
$cell.ClearFormats()
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.3277'
$cell.ClearFormats()
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.0004653'
$cell.ClearFormats()
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.03883'
$cell.ClearFormats()
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.006736'
$cell.ClearFormats()
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.002421'
$cell.ClearFormats()
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.007450'
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.00005596'
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000751'
$cell.ClearFormats()
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.3803'
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.1448'
$cell.ClearFormats()
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.00002102'
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.01011'
$cell.ClearFormats()
